$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.230.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -5.30%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.112.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.17%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'519.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -7.04%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'133.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -6.94%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.117.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -6.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.445"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -6.83%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -9.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -9.85%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -7.65%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.654.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -5.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -2.52%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'25.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -6.56%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.104.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -6.28%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'57.293.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.15%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.0000148"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -10.70%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -7.43%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -11.47%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -8.35%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'343.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -8.44%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'68.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -7.66%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.502"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -7.85%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.236.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -5.80%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.165"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.54%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.21%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0₃0919"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -11.50%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -8.20%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -9.82%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -9.84%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'21.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.16%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.68%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -8.91%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'156.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -6.12%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -9.39%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -11.61%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'25.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.81%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0682"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -8.22%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.147.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.95%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -4.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.680"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -9.78%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.85%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.13%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -6.82%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -9.59%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.252.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.74%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'6.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -6.25%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'19.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -8.38%  "
$ws.Range("E51").Style = "Normal"
